# Generate Report for Archive
# Updates the "Status" values from "Ready for handoff" to "In Translation"
# on the Overview, zh-cn, and de-de sheets, and shrinks the now-narrower
# Status/zh-cn/de-de columns to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# Update the status text wherever it appears.
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"

$ws2.Range("C2").Value = "In Translation"
$ws2.Range("C3").Value = "In Translation"

$ws3.Range("C2").Value = "In Translation"
$ws3.Range("C3").Value = "In Translation"

# Narrow the affected columns to match the shorter status text.
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5

$ws2.Columns.Item(3).ColumnWidth = 12.5

$ws3.Columns.Item(3).ColumnWidth = 12.5
